# Applies the "StoreNo changed to string datatype" update to the Sales
# Reporting workbook: refreshed reporting period/store filter, new
# transaction rows (GZQW2Y8 / WeGift B2B Agency), StoreNo now stored as a
# text code instead of a numeric id, and the recalculated totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Report header / filter summary ------------------------------------
$ws.Range("B2").Value = "GZQW2Y8 | "
$ws.Range("B4").Value = "01-Jul-2024 To 02-Feb-2025"
$ws.Range("B5").Value = 26.65

# ---- Column width tweaks -------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 30.187246322631836
$ws.Columns.Item(8).ColumnWidth = 17.710966110229492

# ---- Row 8 ---------------------------------------------------------------
$ws.Range("A8").Value = 45682.42497685185
$ws.Range("B8").Value = "The Golf Gift Card Variable Digital"
$ws.Range("C8").Value = "EP0130032000173970"
$ws.Range("D8").Formula = "'1435979445"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 1
$ws.Range("G8").Formula = "'GZQW2Y8"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value = "WeGift B2B Agency"
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 6.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 6.5
$ws.Range("O8").Value = 45682.42497685185
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 50
$ws.Range("R8").Value = 0

# ---- Row 9 ---------------------------------------------------------------
$ws.Range("A9").Value = 45681.715208333335
$ws.Range("B9").Value = "The Golf Gift Card Variable Digital"
$ws.Range("C9").Value = "EP0130032000173432"
$ws.Range("D9").Formula = "'1435938493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = 40
$ws.Range("F9").Value = 1
$ws.Range("G9").Formula = "'GZQW2Y8"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value = "WeGift B2B Agency"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 5.2
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 5.2
$ws.Range("O9").Value = 45681.715208333335
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 40
$ws.Range("R9").Value = 0

# ---- Row 10 --------------------------------------------------------------
$ws.Range("A10").Value = 45681.42451388889
$ws.Range("B10").Value = "The Golf Gift Card Variable Digital"
$ws.Range("C10").Value = "EP0130032000173099"
$ws.Range("D10").Formula = "'1435618965"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = 25
$ws.Range("F10").Value = 1
$ws.Range("G10").Formula = "'GZQW2Y8"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value = "WeGift B2B Agency"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 3.25
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 3.25
$ws.Range("O10").Value = 45681.42451388889
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = 0

# ---- Row 11 --------------------------------------------------------------
$ws.Range("A11").Value = 45681.36324074074
$ws.Range("B11").Value = "The Golf Gift Card Variable Digital"
$ws.Range("C11").Value = "EP0130032000173074"
$ws.Range("D11").Formula = "'1435609037"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = 40
$ws.Range("F11").Value = 1
$ws.Range("G11").Formula = "'GZQW2Y8"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value = "WeGift B2B Agency"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 5.2
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 5.2
$ws.Range("O11").Value = 45681.36324074074
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 40
$ws.Range("R11").Value = 0

# ---- Row 12 --------------------------------------------------------------
$ws.Range("A12").Value = 45680.65782407407
$ws.Range("B12").Value = "The Golf Gift Card Variable Digital"
$ws.Range("C12").Value = "EP0130032000172553"
$ws.Range("D12").Formula = "'1435536540"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = 50
$ws.Range("F12").Value = 1
$ws.Range("G12").Formula = "'GZQW2Y8"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value = "WeGift B2B Agency"
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 6.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 6.5
$ws.Range("O12").Value = 45680.65782407407
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = 0

# ---- Row 13 (Totals) ------------------------------------------------------
$ws.Range("E13").Value = 205
$ws.Range("J13").Value = 26.65
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 26.65
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 205
$ws.Range("R13").Value = 0
